$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194; this shifts the existing rows 194-239 down to 195-240
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new record (values copied from the
# "template" of a typical row in this table, with the specific differing fields
# set to the new values from the edit).
$ws.Cells.Item(194, 1).Value = 3
$ws.Cells.Item(194, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(194, 3).Value = 'Coquimbo'
$ws.Cells.Item(194, 4).Value = 44841
$ws.Cells.Item(194, 5).Value = 5
$ws.Cells.Item(194, 6).Value = 'Fruta'
$ws.Cells.Item(194, 7).Value = 100101
$ws.Cells.Item(194, 8).Value = 'Berries'
$ws.Cells.Item(194, 9).Value = 100101001
$ws.Cells.Item(194, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(194, 11).Value = 'Sin especificar'
$ws.Cells.Item(194, 12).Value = 'Primera'
$ws.Cells.Item(194, 13).Value = 40
$ws.Cells.Item(194, 14).Value = 12000
$ws.Cells.Item(194, 15).Value = 12000
$ws.Cells.Item(194, 16).Value = 12000
$ws.Cells.Item(194, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(194, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(194, 19).Value = 6000
$ws.Cells.Item(194, 20).Value = 2
